# Apply trade-value corrections on the "trades_long" sheet and update the
# saved cursor/selection, mirroring the manual edits made in Excel after the
# first run of the app.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected trade parameters (row 2)
$ws.Range("B2").Value = 18000    # Entry price: 20000 -> 18000
$ws.Range("C2").Value = 17900    # Stop loss:   19800 -> 17900
$ws.Range("E2").Value = 0.001    # Position:    0.01  -> 0.001

# Slightly widen the sheet's default column width, as recorded by Excel
# after the edits above.
$ws.StandardWidth = 11.70703125

# Leave the cursor/selection on C2, as saved in the workbook.
$ws.Range("C2").Select() | Out-Null
